$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 105; this shifts the former rows 105..174
# down to 106..175, and the workbook's used range grows from A1:R174 to
# A1:R175 automatically.
$ws.Rows.Item(105).Insert()

# The new row 105 takes on the values that the former row 104 held
# (i.e. the original Albahaca / Región de La Araucanía / $/docena entry
# dated 44210 is duplicated one row down).
$ws.Cells.Item(105,1).Value  = 10
$ws.Cells.Item(105,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(105,3).Value  = "La Araucanía"
$ws.Cells.Item(105,4).Value  = 44210
$ws.Cells.Item(105,5).Value  = 9
$ws.Cells.Item(105,6).Value  = 100112052
$ws.Cells.Item(105,7).Value  = "Albahaca"
$ws.Cells.Item(105,8).Value  = "Sin especificar"
$ws.Cells.Item(105,9).Value  = "Primera"
$ws.Cells.Item(105,10).Value = 100
$ws.Cells.Item(105,11).Value = 5000
$ws.Cells.Item(105,12).Value = 6000
$ws.Cells.Item(105,13).Value = 5450
$ws.Cells.Item(105,14).Value = "$/docena"
$ws.Cells.Item(105,15).Value = "Región de La Araucanía"
$ws.Cells.Item(105,16).Value = 4542
$ws.Cells.Item(105,17).Value = 1.2
$ws.Cells.Item(105,18).Value = "Hortaliza"

# Row 104 (the original row, now sitting above the duplicate) is updated
# with its new reported values: new date, new volume/price figures, unit
# switched from $/docena to $/paquete; the region stays the same.
$ws.Cells.Item(104,4).Value  = 44572
$ws.Cells.Item(104,10).Value = 35
$ws.Cells.Item(104,11).Value = 6000
$ws.Cells.Item(104,12).Value = 6000
$ws.Cells.Item(104,13).Value = 6000
$ws.Cells.Item(104,14).Value = "$/paquete"
$ws.Cells.Item(104,16).Value = 6000
$ws.Cells.Item(104,17).Value = 1
